$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.062304685384816
$ws.Cells.Item(2, 4).Value = 1.064189260177001
$ws.Cells.Item(2, 5).Value = 0.992614727750844
$ws.Cells.Item(2, 6).Value = 1.073736894670149
$ws.Cells.Item(2, 9).Value = 1.052074938943202
$ws.Cells.Item(2, 10).Value = 1.067275896249443
$ws.Cells.Item(2, 11).Value = 1.066905869256302
$ws.Cells.Item(2, 12).Value = 0.9955398523335997
$ws.Cells.Item(2, 13).Value = 1.076427963749218
$ws.Cells.Item(2, 14).Value = 1.068791550440374

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.063422942059306
$ws.Cells.Item(3, 4).Value = 1.065059381542899
$ws.Cells.Item(3, 5).Value = 0.9936372048519299
$ws.Cells.Item(3, 6).Value = 1.074747372887864
$ws.Cells.Item(3, 9).Value = 1.052406041411569
$ws.Cells.Item(3, 10).Value = 1.068047730986211
$ws.Cells.Item(3, 11).Value = 1.067591009589813
$ws.Cells.Item(3, 12).Value = 0.9963617723202687
$ws.Cells.Item(3, 13).Value = 1.077254957915052
$ws.Cells.Item(3, 14).Value = 1.069564481270999

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.064146453422601
$ws.Cells.Item(4, 4).Value = 1.065622352530646
$ws.Cells.Item(4, 5).Value = 0.9942998659930998
$ws.Cells.Item(4, 6).Value = 1.075401476374802
$ws.Cells.Item(4, 9).Value = 1.052619105517238
$ws.Cells.Item(4, 10).Value = 1.06854649309594
$ws.Cells.Item(4, 11).Value = 1.068033643414514
$ws.Cells.Item(4, 12).Value = 0.9968940712668347
$ws.Cells.Item(4, 13).Value = 1.077789719188514
$ws.Cells.Item(4, 14).Value = 1.070063951680132

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.064450600179151
$ws.Cells.Item(5, 4).Value = 1.065859012029514
$ws.Cells.Item(5, 5).Value = 0.994578699834602
$ws.Cells.Item(5, 6).Value = 1.075676522803959
$ws.Cells.Item(5, 9).Value = 1.052708394878483
$ws.Cells.Item(5, 10).Value = 1.068756013269823
$ws.Cells.Item(5, 11).Value = 1.068219559646309
$ws.Cells.Item(5, 12).Value = 0.9971179600053012
$ws.Cells.Item(5, 13).Value = 1.078014446810909
$ws.Cells.Item(5, 14).Value = 1.070273769396694

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.064501666755296
$ws.Cells.Item(6, 4).Value = 1.065898747404022
$ws.Cells.Item(6, 5).Value = 0.994625531979634
$ws.Cells.Item(6, 6).Value = 1.075722707903217
$ws.Cells.Item(6, 9).Value = 1.052723370372996
$ws.Cells.Item(6, 10).Value = 1.068791183285413
$ws.Cells.Item(6, 11).Value = 1.068250766005776
$ws.Cells.Item(6, 12).Value = 0.9971555583673455
$ws.Cells.Item(6, 13).Value = 1.078052174523166
$ws.Cells.Item(6, 14).Value = 1.070308989357741

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.064150517514074
$ws.Cells.Item(7, 4).Value = 1.065625514839722
$ws.Cells.Item(7, 5).Value = 0.994303590798249
$ws.Cells.Item(7, 6).Value = 1.075405151317796
$ws.Cells.Item(7, 9).Value = 1.052620299717184
$ws.Cells.Item(7, 10).Value = 1.068549293340951
$ws.Cells.Item(7, 11).Value = 1.068036128292577
$ws.Cells.Item(7, 12).Value = 0.9968970624462089
$ws.Cells.Item(7, 13).Value = 1.077792722348003
$ws.Cells.Item(7, 14).Value = 1.070066755901812

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.062682621483013
$ws.Cells.Item(8, 4).Value = 1.064483333205005
$ws.Cells.Item(8, 5).Value = 0.9929600610674297
$ws.Cells.Item(8, 6).Value = 1.074078337351246
$ws.Cells.Item(8, 9).Value = 1.052187081189714
$ws.Cells.Item(8, 10).Value = 1.067536879529513
$ws.Cells.Item(8, 11).Value = 1.067137560135085
$ws.Cells.Item(8, 12).Value = 0.9958175282591056
$ws.Cells.Item(8, 13).Value = 1.076707524678425
$ws.Cells.Item(8, 14).Value = 1.069052904346637

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.060095393584969
$ws.Cells.Item(9, 4).Value = 1.062470243236722
$ws.Cells.Item(9, 5).Value = 0.9906006454969559
$ws.Cells.Item(9, 6).Value = 1.071742289412711
$ws.Cells.Item(9, 9).Value = 1.051414645867685
$ws.Cells.Item(9, 10).Value = 1.065747764130038
$ws.Cells.Item(9, 11).Value = 1.065548822530413
$ws.Cells.Item(9, 12).Value = 0.9939188001724441
$ws.Cells.Item(9, 13).Value = 1.074792516269514
$ws.Cells.Item(9, 14).Value = 1.06726124819808

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.058370122510355
$ws.Cells.Item(10, 4).Value = 1.061127907243971
$ws.Cells.Item(10, 5).Value = 0.989033133672735
$ws.Cells.Item(10, 6).Value = 1.070186248409136
$ws.Cells.Item(10, 9).Value = 1.050893601672239
$ws.Cells.Item(10, 10).Value = 1.064551569514875
$ws.Cells.Item(10, 11).Value = 1.064486065401274
$ws.Cells.Item(10, 12).Value = 0.9926553831429383
$ws.Cells.Item(10, 13).Value = 1.073513991204666
$ws.Cells.Item(10, 14).Value = 1.066063354849358

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.05762294308273
$ws.Cells.Item(11, 4).Value = 1.060546595927285
$ws.Cells.Item(11, 5).Value = 0.988355674866747
$ws.Cells.Item(11, 6).Value = 1.06951277884058
$ws.Cells.Item(11, 9).Value = 1.050666538739368
$ws.Cells.Item(11, 10).Value = 1.064032781712134
$ws.Cells.Item(11, 11).Value = 1.064025024458291
$ws.Cells.Item(11, 12).Value = 0.9921088820399291
$ws.Cells.Item(11, 13).Value = 1.072959934749653
$ws.Cells.Item(11, 14).Value = 1.065543830308432

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.057345387464081
$ws.Cells.Item(12, 4).Value = 1.060330660288808
$ws.Cells.Item(12, 5).Value = 0.9881042295826724
$ws.Cells.Item(12, 6).Value = 1.069262668120909
$ws.Cells.Item(12, 9).Value = 1.050581979834142
$ws.Cells.Item(12, 10).Value = 1.063839956050347
$ws.Cells.Item(12, 11).Value = 1.06385364385193
$ws.Cells.Item(12, 12).Value = 0.9919059725120875
$ws.Cells.Item(12, 13).Value = 1.072754066168675
$ws.Cells.Item(12, 14).Value = 1.065350730812088

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.057404925004029
$ws.Cells.Item(13, 4).Value = 1.060376979703082
$ws.Cells.Item(13, 5).Value = 0.9881581567098651
$ws.Cells.Item(13, 6).Value = 1.069316315646699
$ws.Cells.Item(13, 9).Value = 1.050600127856974
$ws.Cells.Item(13, 10).Value = 1.063881323462489
$ws.Cells.Item(13, 11).Value = 1.063890411436722
$ws.Cells.Item(13, 12).Value = 0.9919494934313052
$ws.Cells.Item(13, 13).Value = 1.072798228729428
$ws.Cells.Item(13, 14).Value = 1.065392156970699

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.057600000660845
$ws.Cells.Item(14, 4).Value = 1.060528746839395
$ws.Cells.Item(14, 5).Value = 0.9883348863814464
$ws.Cells.Item(14, 6).Value = 1.069492103675877
$ws.Cells.Item(14, 9).Value = 1.050659553514301
$ws.Cells.Item(14, 10).Value = 1.064016845230587
$ws.Cells.Item(14, 11).Value = 1.064010860726565
$ws.Cells.Item(14, 12).Value = 0.9920921077337197
$ws.Cells.Item(14, 13).Value = 1.072942918963706
$ws.Cells.Item(14, 14).Value = 1.065527871195253

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.057720190542096
$ws.Cells.Item(15, 4).Value = 1.060622254156748
$ws.Cells.Item(15, 5).Value = 0.9884438009545853
$ws.Cells.Item(15, 6).Value = 1.06960041855455
$ws.Cells.Item(15, 9).Value = 1.050696138780901
$ws.Cells.Item(15, 10).Value = 1.064100328118252
$ws.Cells.Item(15, 11).Value = 1.064085056345459
$ws.Cells.Item(15, 12).Value = 0.9921799884222134
$ws.Cells.Item(15, 13).Value = 1.073032058461541
$ws.Cells.Item(15, 14).Value = 1.065611472638193

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.058419707612497
$ws.Cells.Item(16, 4).Value = 1.061166485460325
$ws.Cells.Item(16, 5).Value = 0.9890781214508737
$ws.Cells.Item(16, 6).Value = 1.070230950807944
$ws.Cells.Item(16, 9).Value = 1.050908640582337
$ws.Cells.Item(16, 10).Value = 1.064585982278687
$ws.Cells.Item(16, 11).Value = 1.064516645027498
$ws.Cells.Item(16, 12).Value = 0.9926916645766087
$ws.Cells.Item(16, 13).Value = 1.073550752675839
$ws.Cells.Item(16, 14).Value = 1.066097816483242

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.058858461817303
$ws.Cells.Item(17, 4).Value = 1.061507848328666
$ws.Cells.Item(17, 5).Value = 0.989476357848556
$ws.Cells.Item(17, 6).Value = 1.070626548954285
$ws.Cells.Item(17, 9).Value = 1.05104154972718
$ws.Cells.Item(17, 10).Value = 1.064890398417138
$ws.Cells.Item(17, 11).Value = 1.064787138755749
$ws.Cells.Item(17, 12).Value = 0.9930127773699352
$ws.Cells.Item(17, 13).Value = 1.073875996100402
$ws.Cells.Item(17, 14).Value = 1.066402664927526

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.059114367665895
$ws.Cells.Item(18, 4).Value = 1.06170695267123
$ws.Cells.Item(18, 5).Value = 0.9897087662937556
$ws.Cells.Item(18, 6).Value = 1.070857324240427
$ws.Cells.Item(18, 9).Value = 1.051118933745224
$ws.Cells.Item(18, 10).Value = 1.065067879373874
$ws.Cells.Item(18, 11).Value = 1.06494483020157
$ws.Cells.Item(18, 12).Value = 0.9932001317071769
$ws.Cells.Item(18, 13).Value = 1.074065662006596
$ws.Cells.Item(18, 14).Value = 1.066580397927578

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.059201622933536
$ws.Cells.Item(19, 4).Value = 1.061774840948091
$ws.Cells.Item(19, 5).Value = 0.9897880325774034
$ws.Cells.Item(19, 6).Value = 1.070936017671931
$ws.Cells.Item(19, 9).Value = 1.051145296015921
$ws.Cells.Item(19, 10).Value = 1.065128382266687
$ws.Cells.Item(19, 11).Value = 1.064998584844105
$ws.Cells.Item(19, 12).Value = 0.9932640239640975
$ws.Cells.Item(19, 13).Value = 1.07413032587919
$ws.Cells.Item(19, 14).Value = 1.066640986741439

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.058811388897716
$ws.Cells.Item(20, 4).Value = 1.06147122403172
$ws.Cells.Item(20, 5).Value = 0.9894336180360679
$ws.Cells.Item(20, 6).Value = 1.070584101958594
$ws.Cells.Item(20, 9).Value = 1.051027304278322
$ws.Cells.Item(20, 10).Value = 1.064857745706247
$ws.Cells.Item(20, 11).Value = 1.064758125922272
$ws.Cells.Item(20, 12).Value = 0.9929783193494215
$ws.Cells.Item(20, 13).Value = 1.073841105016649
$ws.Cells.Item(20, 14).Value = 1.06636996584604

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.057542556299963
$ws.Cells.Item(21, 4).Value = 1.060484055497287
$ws.Cells.Item(21, 5).Value = 0.9882828385668249
$ws.Cells.Item(21, 6).Value = 1.069440337218033
$ws.Cells.Item(21, 9).Value = 1.050642060148179
$ws.Cells.Item(21, 10).Value = 1.063976940895462
$ws.Cells.Item(21, 11).Value = 1.06397539498893
$ws.Cells.Item(21, 12).Value = 0.9920501090198102
$ws.Cells.Item(21, 13).Value = 1.072900313156965
$ws.Cells.Item(21, 14).Value = 1.065487910191395

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.056744675177161
$ws.Cells.Item(22, 4).Value = 1.059863321328022
$ws.Cells.Item(22, 5).Value = 0.9875604150241495
$ws.Cells.Item(22, 6).Value = 1.068721472198179
$ws.Cells.Item(22, 9).Value = 1.050398582694799
$ws.Cells.Item(22, 10).Value = 1.063422421842346
$ws.Cells.Item(22, 11).Value = 1.063482511772369
$ws.Cells.Item(22, 12).Value = 0.9914670000341481
$ws.Cells.Item(22, 13).Value = 1.072308410223782
$ws.Cells.Item(22, 14).Value = 1.064932603657619

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.057167658368698
$ws.Cells.Item(23, 4).Value = 1.060192390168901
$ws.Cells.Item(23, 5).Value = 0.9879432794643023
$ws.Cells.Item(23, 6).Value = 1.069102531061668
$ws.Cells.Item(23, 9).Value = 1.050527774112634
$ws.Cells.Item(23, 10).Value = 1.063716451475291
$ws.Cells.Item(23, 11).Value = 1.063743869592464
$ws.Cells.Item(23, 12).Value = 0.991776070289318
$ws.Cells.Item(23, 13).Value = 1.0726222261914
$ws.Cells.Item(23, 14).Value = 1.065227050846369

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.058832659160757
$ws.Cells.Item(24, 4).Value = 1.061487772996039
$ws.Cells.Item(24, 5).Value = 0.9894529299347244
$ws.Cells.Item(24, 6).Value = 1.070603281837269
$ws.Cells.Item(24, 9).Value = 1.051033741614826
$ws.Cells.Item(24, 10).Value = 1.064872500307226
$ws.Cells.Item(24, 11).Value = 1.064771235829078
$ws.Cells.Item(24, 12).Value = 0.9929938892766442
$ws.Cells.Item(24, 13).Value = 1.07385687092847
$ws.Cells.Item(24, 14).Value = 1.066384741400246

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.060764330542578
$ws.Cells.Item(25, 4).Value = 1.062990724353379
$ws.Cells.Item(25, 5).Value = 0.9912096547607049
$ws.Cells.Item(25, 6).Value = 1.072345980469564
$ws.Cells.Item(25, 9).Value = 1.051615411602341
$ws.Cells.Item(25, 10).Value = 1.066210900627256
$ws.Cells.Item(25, 11).Value = 1.065960183243601
$ws.Cells.Item(25, 12).Value = 0.9944092447426414
$ws.Cells.Item(25, 13).Value = 1.075287918198007
$ws.Cells.Item(25, 14).Value = 1.067725042402246
